$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: headers (reordered / relabeled) ---
# Set in this particular order so newly-introduced shared strings land
# in the same sequence as the authored workbook.
$ws.Range("F1").Value = "CPL_A"
$ws.Range("E1").Value = "Z_D"
$ws.Range("D1").Value = "CPL_D"
$ws.Range("C1").Value = "AND_A_D"
$ws.Range("B1").Value = "CPL_OUT"
$ws.Range("A1").Value = "CATCH_FLAGS"
$ws.Range("G1").Value = "Z_A"
$ws.Range("H1").Value = "f(x,y)"

# --- Row 2: descriptions (reordered) ---
$ws.Range("A2").Value = "if set flags will be saved as long as"
$ws.Range("B2").Value = "complement the output"
$ws.Range("C2").Value = "if set X&Y else X+Y"
$ws.Range("D2").Value = "complement the input data"
$ws.Range("E2").Value = "If set ignore what input on databus and input zero"
$ws.Range("F2").Value = "complement A register input to the ALU after Za"
$ws.Range("G2").Value = "ignore A register and make its input Zero"

# --- Selection moves from C15 to B7 ---
$ws.Range("B7").Select()
